$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8709.154
$ws.Range("I62").Value = 1136.9
$ws.Range("J62").Value = 33950
$ws.Range("K62").Value = 1136.9
$ws.Range("L62").Value = 33950
$ws.Range("M62").Value = -512.9000000000001
$ws.Range("N62").Value = -35198
$ws.Range("H65").Value = 8709.154
$ws.Range("I65").Value = 1136.9
$ws.Range("J65").Value = 33950
$ws.Range("K65").Value = 5684.5
$ws.Range("L65").Value = 169750
$ws.Range("M65").Value = -2564.5
$ws.Range("N65").Value = -175990
$ws.Range("H96").Value = 635.1111
$ws.Range("I96").Value = 500
$ws.Range("J96").Value = 673.7143
$ws.Range("K96").Value = 1500
$ws.Range("L96").Value = 2021.1429
$ws.Range("M96").Value = -127
$ws.Range("N96").Value = -4767.1429
$ws.Range("H99").Value = 359.75
$ws.Range("I99").Value = 219.5
$ws.Range("K99").Value = 658.5
$ws.Range("M99").Value = 839.5
$ws.Range("H101").Value = 2778.25
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H137").Value = 1539.1111
$ws.Range("I137").Value = 1206.1613
$ws.Range("J137").Value = 1987.8695
$ws.Range("K137").Value = 3618.4839
$ws.Range("L137").Value = 5963.6085
$ws.Range("M137").Value = -1068.4839
$ws.Range("N137").Value = -11063.6085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 22218.777
$ws.Range("J37").Value = 26994.834
$ws.Range("L37").Value = 26994.834
$ws.Range("N37").Value = -27540.834
$ws.Range("H45").Value = 2886.0605
$ws.Range("I45").Value = 1993
$ws.Range("K45").Value = 1993
$ws.Range("M45").Value = -1616
$ws.Range("H80").Value = 20828.6
$ws.Range("J80").Value = 20828.6
$ws.Range("L80").Value = 20828.6
$ws.Range("N80").Value = -22824.6
$ws.Range("H83").Value = 20828.6
$ws.Range("J83").Value = 20828.6
$ws.Range("L83").Value = 62485.8
$ws.Range("N83").Value = -72469.79999999999
$ws.Range("H117").Value = 30537.2
$ws.Range("J117").Value = 30537.2
$ws.Range("L117").Value = 30537.2
$ws.Range("N117").Value = -39715.2
$ws.Range("H124").Value = 22256.143
$ws.Range("J124").Value = 22256.143
$ws.Range("L124").Value = 22256.143
$ws.Range("N124").Value = -32076.143
$ws.Range("H125").Value = 900000000
$ws.Range("J125").Value = 900000000
$ws.Range("L125").Value = 900000000
$ws.Range("N125").Value = -900009840
$ws.Range("H132").Value = 2667.6191
$ws.Range("I132").Value = 2502.923
$ws.Range("J132").Value = 2935.25
$ws.Range("K132").Value = 7508.768999999999
$ws.Range("L132").Value = 8805.75
$ws.Range("M132").Value = -4978.768999999999
$ws.Range("N132").Value = -13865.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 10030000
$ws.Range("J9").Value = 10030000
$ws.Range("L9").Value = 10030000
$ws.Range("N9").Value = -10030336
$ws.Range("H44").Value = 17250
$ws.Range("J44").Value = 17250
$ws.Range("L44").Value = 17250
$ws.Range("N44").Value = -18244
$ws.Range("H82").Value = 23709
$ws.Range("I82").Value = 7338.6665
$ws.Range("J82").Value = 29165.777
$ws.Range("K82").Value = 7338.6665
$ws.Range("L82").Value = 29165.777
$ws.Range("M82").Value = -6955.6665
$ws.Range("N82").Value = -29931.777
$ws.Range("H85").Value = 23709
$ws.Range("I85").Value = 7338.6665
$ws.Range("J85").Value = 29165.777
$ws.Range("K85").Value = 7338.6665
$ws.Range("L85").Value = 29165.777
$ws.Range("M85").Value = -6012.6665
$ws.Range("N85").Value = -31817.777

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 3493
$ws.Range("I35").Value = 1604.125
$ws.Range("J35").Value = 8530
$ws.Range("K35").Value = 1604.125
$ws.Range("L35").Value = 8530
$ws.Range("M35").Value = -1310.125
$ws.Range("N35").Value = -9118
$ws.Range("H38").Value = 4361.3335
$ws.Range("J38").Value = 4361.3335
$ws.Range("L38").Value = 4361.3335
$ws.Range("N38").Value = -5115.3335
$ws.Range("H41").Value = 15909
$ws.Range("J41").Value = 19761.25
$ws.Range("L41").Value = 19761.25
$ws.Range("N41").Value = -20617.25
$ws.Range("H46").Value = 4361.3335
$ws.Range("J46").Value = 4361.3335
$ws.Range("L46").Value = 4361.3335
$ws.Range("N46").Value = -4783.3335
$ws.Range("H50").Value = 8906.571
$ws.Range("J50").Value = 8906.571
$ws.Range("L50").Value = 8906.571
$ws.Range("N50").Value = -10156.571
$ws.Range("H51").Value = 9490.333000000001
$ws.Range("J51").Value = 9490.333000000001
$ws.Range("L51").Value = 9490.333000000001
$ws.Range("N51").Value = -10962.333
$ws.Range("H60").Value = 33375.5
$ws.Range("J60").Value = 33375.5
$ws.Range("L60").Value = 33375.5
$ws.Range("N60").Value = -34397.5
$ws.Range("H61").Value = 9490.333000000001
$ws.Range("J61").Value = 9490.333000000001
$ws.Range("L61").Value = 9490.333000000001
$ws.Range("N61").Value = -10186.333
$ws.Range("H109").Value = 10925
$ws.Range("J109").Value = 10925
$ws.Range("L109").Value = 10925
$ws.Range("N109").Value = -13005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 34500
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H97").Value = 2822.5
$ws.Range("I97").Value = 3296.4285
$ws.Range("J97").Value = 1716.6666
$ws.Range("K97").Value = 3296.4285
$ws.Range("L97").Value = 1716.6666
$ws.Range("M97").Value = -2800.4285
$ws.Range("N97").Value = -2708.6666
$ws.Range("H122").Value = 2112
$ws.Range("J122").Value = 2200
$ws.Range("L122").Value = 6600
$ws.Range("N122").Value = -11500
$ws.Range("H123").Value = 22873.2
$ws.Range("J123").Value = 22873.2
$ws.Range("L123").Value = 22873.2
$ws.Range("N123").Value = -27773.2
$ws.Range("H136").Value = 27611.615
$ws.Range("J136").Value = 27611.615
$ws.Range("L136").Value = 82834.845
$ws.Range("N136").Value = -87934.845
$ws.Range("H139").Value = 30249.5
$ws.Range("J139").Value = 30249.5
$ws.Range("L139").Value = 30249.5
$ws.Range("N139").Value = -40529.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5884465.5
$ws.Range("I7").Value = 8335076
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 8335076
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -8334964
$ws.Range("N7").Value = -3224
$ws.Range("H22").Value = 1345.5454
$ws.Range("I22").Value = 1000.1667
$ws.Range("J22").Value = 1760
$ws.Range("K22").Value = 1000.1667
$ws.Range("L22").Value = 1760
$ws.Range("M22").Value = -705.1667
$ws.Range("N22").Value = -2350
$ws.Range("H27").Value = 1345.5454
$ws.Range("I27").Value = 1000.1667
$ws.Range("J27").Value = 1760
$ws.Range("K27").Value = 1000.1667
$ws.Range("L27").Value = 1760
$ws.Range("M27").Value = -893.1667
$ws.Range("N27").Value = -1974
$ws.Range("H61").Value = 2468.0588
$ws.Range("I61").Value = 1795.7
$ws.Range("K61").Value = 1795.7
$ws.Range("M61").Value = -1593.7
$ws.Range("H107").Value = 2250
$ws.Range("I107").Value = 2250
$ws.Range("K107").Value = 2250
$ws.Range("M107").Value = -330
$ws.Range("H113").Value = 2468.0588
$ws.Range("I113").Value = 1795.7
$ws.Range("K113").Value = 1795.7
$ws.Range("M113").Value = 374.3
$ws.Range("H126").Value = 5884465.5
$ws.Range("I126").Value = 8335076
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 25005228
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -25002758
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 47465.6
$ws.Range("J46").Value = 47465.6
$ws.Range("L46").Value = 47465.6
$ws.Range("N46").Value = -47927.6
$ws.Range("H62").Value = 6224.375
$ws.Range("I62").Value = 3662.7273
$ws.Range("J62").Value = 11860
$ws.Range("K62").Value = 3662.7273
$ws.Range("L62").Value = 11860
$ws.Range("M62").Value = -3038.7273
$ws.Range("N62").Value = -13108
$ws.Range("H65").Value = 6224.375
$ws.Range("I65").Value = 3662.7273
$ws.Range("J65").Value = 11860
$ws.Range("K65").Value = 18313.6365
$ws.Range("L65").Value = 59300
$ws.Range("M65").Value = -15193.6365
$ws.Range("N65").Value = -65540
$ws.Range("H109").Value = 24600
$ws.Range("J109").Value = 24600
$ws.Range("L109").Value = 24600
$ws.Range("N109").Value = -27374
$ws.Range("H119").Value = 29666
$ws.Range("J119").Value = 29666
$ws.Range("L119").Value = 29666
$ws.Range("N119").Value = -39342
$ws.Range("H123").Value = 18409.666
$ws.Range("J123").Value = 18409.666
$ws.Range("L123").Value = 18409.666
$ws.Range("N123").Value = -28209.666
$ws.Range("H125").Value = 55992.5
$ws.Range("J125").Value = 55992.5
$ws.Range("L125").Value = 55992.5
$ws.Range("N125").Value = -65832.5
$ws.Range("H134").Value = 47465.6
$ws.Range("J134").Value = 47465.6
$ws.Range("L134").Value = 142396.8
$ws.Range("N134").Value = -147466.8
